# Update the three patient records on Sheet1 (A=id, B=name, M=timestamp,
# P=gender, Q=age) to the corrected values captured after fixing the
# submission/edit form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 (Ahmed / male): only the submission timestamp changed.
$ws.Range("M1").NumberFormat = "d/m/yy hh:mm"
$ws.Range("M1").Value = 42457.078729479166
$ws.Range("Q1").Value = 77

# Row 2: name corrected to "msel7y", gender corrected female -> male,
# age corrected 19 -> 90, timestamp refreshed.
$ws.Range("B2").Value = "msel7y"
$ws.Range("M2").NumberFormat = "d/m/yy hh:mm"
$ws.Range("M2").Value = 42457.07872990741
$ws.Range("P2").Value = "male"
$ws.Range("Q2").Value = 90

# Row 3: name corrected to "ayaaaa", gender corrected male -> female,
# age corrected 18 -> 22, timestamp refreshed.
$ws.Range("B3").Value = "ayaaaa"
$ws.Range("M3").NumberFormat = "d/m/yy hh:mm"
$ws.Range("M3").Value = 42457.07873026621
$ws.Range("P3").Value = "female"
$ws.Range("Q3").Value = 22
